$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 82: date serial 45884 in column A (formatted like the row above it)
# and the numeric value in column B.
$ws.Range("A82").Value = 45884
$ws.Range("B82").Value = 0.06654624964350926

# Copy the style (number format, font, border, alignment) of A81 onto A82
$ws.Range("A81").Copy()
$ws.Range("A82").PasteSpecial(-4122) # xlPasteFormats
